$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2029")

# Fill in the previously-blank dev_trait (column J) cells for the rows
# that already had a K (overall_start) value but no J value yet.
$ws.Range("J5").Value = "star"
$ws.Range("J6").Value = "elite"
$ws.Range("J8").Value = "elite"
$ws.Range("J14").Value = "impact"
$ws.Range("J15").Value = "impact"
$ws.Range("J16").Value = "elite"
$ws.Range("J17").Value = "star"
$ws.Range("J24").Value = "elite"
$ws.Range("J25").Value = "impact"
$ws.Range("J26").Value = "normal"
$ws.Range("J28").Value = "star"
$ws.Range("J29").Value = "impact"
$ws.Range("J34").Value = "elite"
$ws.Range("J36").Value = "star"
$ws.Range("J43").Value = "impact"
$ws.Range("J46").Value = "impact"
$ws.Range("J51").Value = "normal"
$ws.Range("J55").Value = "normal"
$ws.Range("J59").Value = "normal"
$ws.Range("J63").Value = "impact"
$ws.Range("J68").Value = "star"
$ws.Range("J70").Value = "star"
$ws.Range("J71").Value = "star"
$ws.Range("J77").Value = "impact"
$ws.Range("J82").Value = "star"
$ws.Range("J89").Value = "elite"
$ws.Range("J90").Value = "impact"
$ws.Range("J92").Value = "star"
$ws.Range("J97").Value = "normal"
$ws.Range("J101").Value = "star"
$ws.Range("J103").Value = "impact"

# Move the active selection to where editing left off (matches the
# author's final cursor position in the saved workbook).
$ws.Range("J105").Select()
